$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: Average of J column (bold, size 11)
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"
$ws.Range("J12").Font.Bold = $true

# Row 14-17: labels and summary formulas (bold, size 12, vertical centered)
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

$ws.Range("B14:B17").Font.Bold = $true
$ws.Range("B14:B17").Font.Size = 12
$ws.Range("B14:B17").VerticalAlignment = -4108
$ws.Range("A14:B17").RowHeight = 15.6

$ws.Range("A14").Select()
$ws.Range("A14:B17").Select()

$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
